# Apply the "Fixed update to excel issue" change:
#  1. Rename "Requested quantity" header -> "Weekly_PO_Qty" on "Weekly Quantity" sheet
#  2. Rename "Requested quantity" header -> "Monthly_PO_Qty" on "Monthly Trend" sheet
#  3. Add a new "PO Forecast" sheet (as the last sheet) with forecast data

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity sheet header rename ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend sheet header rename ---
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. Add new "PO Forecast" sheet at the end ---
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$data = @(
    @(44934.99999999999, 273, -109.9846751838848, 648.9881601250764),
    @(44962.99999999999, 335, -50.19628163452771, 752.99637847989),
    @(44976.99999999999, 366, -31.39415455410789, 756.6660287376294),
    @(44990.99999999999, 397, 16.98171660226523, 753.0732887244438),
    @(44997.99999999999, 413, 46.88437717973782, 802.4806902991605),
    @(45011.99999999999, 444, 62.37350457283191, 815.7612450105506),
    @(45060.99999999999, 553, 143.1101479200495, 954.7901555848266),
    @(45130.99999999999, 709, 327.7549661613439, 1108.924183741586),
    @(45137.99999999999, 725, 335.712644160598, 1111.501044607656),
    @(45144.99999999999, 740, 343.3156079079068, 1126.318289429959),
    @(45151.99999999999, 756, 347.9317609615703, 1147.901285628977),
    @(45158.99999999999, 771, 390.7571247385469, 1177.742180431172),
    @(45165.99999999999, 787, 395.3023308094364, 1173.228021914204),
    @(45172.99999999999, 803, 407.5639998728102, 1223.939095749454),
    @(45179.99999999999, 818, 426.496551360041, 1211.216734029046),
    @(45186.99999999999, 834, 464.5437736336414, 1227.349577349068)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

# --- Formatting: reuse the existing header / date styles from the Weekly Quantity sheet ---
# Header style (bold, centered, thin border) -> row 1, all 4 columns
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Date style -> column A, data rows
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A17").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
